$d = $word.ActiveDocument

# Wrapper around Find.Execute for a literal (non-wildcard) find-and-replace.
# wdReplace = 2 (wdReplaceAll) so every matching occurrence found by this call is replaced.
function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Remove company "4 ÉP-SZAK 2000 Kft." identity details in the "székhelyéül" sentence
Replace-Text "a(z) 4 ÉP-SZAK 2000 Kft. (1139 Budapest, Országbíró u 2. 4. em. 19., cégj.: 01-09-687257) székhelyéül" "a(z)  (, cégj.: ) székhelyéül"

# 2. Remove company name from the "az ingatlan felett..." sentence
Replace-Text "A 4 ÉP-SZAK 2000 Kft. az ingatlan felett" "A  az ingatlan felett"

# 3. Remove company + representative details in the two "amely létrejött egyrészről..." paragraphs
Replace-Text "amely létrejött egyrészről 4 ÉP-SZAK 2000 Kft. (1139 Budapest, Országbíró u 2. 4. em. 19., cégj.: 01-09-687257), képviseli: Bazsika István, ügyvezetõ (1012 Budapest, Logodi u. 48. fszt. 1., ig.sz.: 457361HA, an.neve: Süle Mária Margit), mint Bérlő," "amely létrejött egyrészről  (, cégj.: ), képviseli: , manager (, ig.sz.: , an.neve: ), mint Bérlő,"
Replace-Text "amely létrejött egyrészről 4 ÉP-SZAK 2000 Kft. (1139 Budapest, Országbíró u 2. 4. em. 19., cégj.: 01-09-687257), képviseli: Bazsika István, ügyvezetõ (1012 Budapest, Logodi u. 48. fszt. 1., ig.sz.: 457361HA, an.neve: Süle Mária Margit), mint Bérlő," "amely létrejött egyrészről  (, cégj.: ), képviseli: , manager (, ig.sz.: , an.neve: ), mint Bérlő,"

# 4. Remove the "cégiratok őrzési címe" address (point 11 of megbízási szerződés)
Replace-Text "A cég iratainak őrzési címe: Bazsika István, 1012 Budapest, Logodi u. 48. fszt. 1." "A cég iratainak őrzési címe: , "

# 5. Fee / term details in "MEGBÍZÁSI SZERZŐDÉS DÍJSZABÁS/FUTAMIDŐ"
Replace-Text "1. A megbízási/cégképviseleti díj nettó 32940 Ft, azaz Harminckettõezerkilencszáznegyven Forint." "1. A megbízási/cégképviseleti díj nettó 0 Ft, azaz  Forint."
Replace-Text "Az induló dátum: 2015-01-25" "Az induló dátum: 2015-10-23"
Replace-Text "A következő díj esedékessége: 2015-07-25" "A következő díj esedékessége: 2015-10-23"
Replace-Text "4. 1 (azaz egy) havi megbízási/cégképviseleti díj összege a szerződés aláírásakor nettó 5490 Forint." "4. 1 (azaz egy) havi megbízási/cégképviseleti díj összege a szerződés aláírásakor nettó 0 Forint."

# 6. Contact details in the final annex (küldemények kezelési rendje)
Replace-Text "Telefon: 06-20/988-93-56" "Telefon: "
Replace-Text "Email cím: bazsika.istvan@gmail.com" "Email cím: "
Replace-Text "Postacím: Bazsika István, 1012 Budapest, Logodi u. 48. fszt. 1." "Postacím: , "

# 7. Update every "Budapest, 2015-01-25" signature-block date stamp (6 occurrences) in one pass
Replace-Text "Budapest, 2015-01-25" "Budapest, 2015-10-23"

Write-Output "done"
